# Feedback Form scenarios : Contact Name, Editing Contact Name Field
$wb = $excel.ActiveWorkbook

# --- Login sheet: C3 value is re-keyed (string content stays "sdf") ---
$loginWs = $wb.Worksheets.Item("Login")
$loginWs.Range("C3").Value = "sdf"

# --- Rename Sheet1 -> Feedback and build out the feedback form data ---
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Name = "Feedback"

# Row 1 headers (A1:H1, then J1:L1, leaving I1 for later)
$ws.Range("A1").Value = "Contact name"
$ws.Range("B1").Value = "Primary role"
$ws.Range("C1").Value = "Recognition level"
$ws.Range("D1").Value = "contact attributes"
$ws.Range("E1").Value = "institution"
$ws.Range("F1").Value = "institution type"
$ws.Range("G1").Value = "Payer"
$ws.Range("H1").Value = "IDS Member"
$ws.Range("J1").Value = "date"
$ws.Range("K1").Value = "time spent"
$ws.Range("L1").Value = "What topics were discussed?"

# E1 gets wrap text alignment
$ws.Range("E1").WrapText = $true

# Row 2 data
$ws.Range("A2").Value = "Nirmalkumar Pant"
$ws.Range("B2").Value = "MCO Manager"
$ws.Range("C2").Value = "National"
$ws.Range("D2").Value = "Cardiology"
$ws.Range("I2").Value = "Adventist Health System"
$ws.Range("F2").Value = "Academic"
$ws.Range("E2").Value = "Florida Hospital"

$ws.Range("G2").Value = $true
$ws.Range("H2").Value = $true

$ws.Range("J2").Value = 42857
$ws.Range("J2").NumberFormat = "mm-dd-yy"

$ws.Range("K2").Value = "10 minutes"

# I1 header filled in last
$ws.Range("I1").Value = "IDS Input"

# Column widths (bestFit, matching the authored layout - closest
# representable values given the host's column-width rounding grid)
$ws.Columns.Item(1).ColumnWidth = 16.5
$ws.Columns.Item(2).ColumnWidth = 11
$ws.Columns.Item(3).ColumnWidth = 15.666666666666666
$ws.Columns.Item(4).ColumnWidth = 15.833333333333334
$ws.Columns.Item(5).ColumnWidth = 22.333333333333332
$ws.Columns.Item(6).ColumnWidth = 14
$ws.Columns.Item(8).ColumnWidth = 11
$ws.Columns.Item(9).ColumnWidth = 22.333333333333332
$ws.Columns.Item(11).ColumnWidth = 9.666666666666666
$ws.Columns.Item(12).ColumnWidth = 26.166666666666668

# Select A2 and make this the active sheet/tab
$ws.Range("A2").Select() | Out-Null
$ws.Activate() | Out-Null
